# This workbook holds a weekly-updated price log for "Betarraga" (Beet)
# at "Vega Monumental Concepción". Each week's refresh prepends one new
# pair of records (quality "Primera" and "Segunda") at the top of the
# data block (rows 26-27), pushing every existing record down by one
# pair-slot (2 rows). The oldest pair, previously at the bottom of the
# range, ends up in the two newly created rows at the end of the sheet.
#
# Concretely: insert two blank rows at row 26 (shifting rows 26:214 -> 28:216,
# which also grows the sheet by the required two rows), then populate the
# freshly inserted rows 26 and 27 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data block (rows 26-214) down by two rows, inserting
# two new blank rows at the top of the block and extending the sheet.
$ws.Range("A26:R27").Insert()

# New "Primera" record for the latest week.
$ws.Cells.Item(26, 1).Value2 = 11
$ws.Cells.Item(26, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(26, 3).Value2 = "Bíobío"
$ws.Cells.Item(26, 4).Value2 = 44490
$ws.Cells.Item(26, 5).Value2 = 8
$ws.Cells.Item(26, 6).Value2 = 100114014
$ws.Cells.Item(26, 7).Value2 = "Betarraga"
$ws.Cells.Item(26, 8).Value2 = "Sin especificar"
$ws.Cells.Item(26, 9).Value2 = "Primera"
$ws.Cells.Item(26, 10).Value2 = 200
$ws.Cells.Item(26, 11).Value2 = 600
$ws.Cells.Item(26, 12).Value2 = 700
$ws.Cells.Item(26, 13).Value2 = 650
$ws.Cells.Item(26, 14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(26, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(26, 16).Value2 = 130
$ws.Cells.Item(26, 17).Value2 = 5
$ws.Cells.Item(26, 18).Value2 = "Hortaliza"

# New "Segunda" record for the latest week.
$ws.Cells.Item(27, 1).Value2 = 11
$ws.Cells.Item(27, 2).Value2 = "Vega Monumental Concepción"
$ws.Cells.Item(27, 3).Value2 = "Bíobío"
$ws.Cells.Item(27, 4).Value2 = 44490
$ws.Cells.Item(27, 5).Value2 = 8
$ws.Cells.Item(27, 6).Value2 = 100114014
$ws.Cells.Item(27, 7).Value2 = "Betarraga"
$ws.Cells.Item(27, 8).Value2 = "Sin especificar"
$ws.Cells.Item(27, 9).Value2 = "Segunda"
$ws.Cells.Item(27, 10).Value2 = 100
$ws.Cells.Item(27, 11).Value2 = 500
$ws.Cells.Item(27, 12).Value2 = 500
$ws.Cells.Item(27, 13).Value2 = 500
$ws.Cells.Item(27, 14).Value2 = "$/paquete 5 unidades"
$ws.Cells.Item(27, 15).Value2 = "Región Metropolitana"
$ws.Cells.Item(27, 16).Value2 = 100
$ws.Cells.Item(27, 17).Value2 = 5
$ws.Cells.Item(27, 18).Value2 = "Hortaliza"
